$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.006.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.55%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.823.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.56%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.40%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.26%  "

$ws.Range("E7").Value = "  +1.22%  "

$ws.Range("E8").Value = "  -0.85%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07356"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8746"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.54%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.30"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.53%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.841.15"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.99%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07300"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.18%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.430"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.63%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.519"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.14%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.80"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.10%  "

$ws.Range("E17").Value = "  +0.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008746"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.46%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.004"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.20%  "

$ws.Range("E20").Value = "  +0.11%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.015.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.291"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.83%  "

$ws.Range("E23").Value = "  +0.70%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.054.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.894"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.13%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.43%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.141"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.88%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.247"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.93%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.36%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08889"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.19%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7550"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.161"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.510"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.19%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.931"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.14%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.005"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.39%  "

$ws.Range("E37").Value = "  +0.24%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05315"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.22%  "

$ws.Range("E39").Value = "  -0.17%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.984"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.44%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.205"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.22%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.370"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.66%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5304"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.45%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1654"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.37%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.479"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.69%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4896"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.76%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.004"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.24%  "

$ws.Range("E49").Value = "  -0.51%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.47%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06301"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.57%  "
